# Update "想去人数" (want-to-go count) values in column F across sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        2 = 1644
        4 = 9448
        5 = 729
        8 = 321
        11 = 1620
        12 = 1396
        15 = 1463
        21 = 373
        22 = 1104
        33 = 162
        38 = 225
        39 = 600
        41 = 1252
    }
    "演出" = @{
        22 = 653
        25 = 302
    }
    "本地生活" = @{
        7 = 2351
        8 = 3550
        9 = 4
        11 = 53
        12 = 89
    }
    "全部类型" = @{
        2 = 1644
        3 = 9448
        6 = 3550
        7 = 729
        8 = 53
        9 = 53
        12 = 321
        14 = 1396
        16 = 89
        17 = 89
        18 = 1463
        22 = 1104
        38 = 162
        39 = 302
        40 = 302
        43 = 225
        45 = 600
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsMap = $updates[$sheetName]
    foreach ($row in $rowsMap.Keys) {
        $ws.Range("F$row").Value = $rowsMap[$row]
    }
}
